$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Adjust scroll position (removes stale horizontal scroll state) ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 1

# --- Replace the "Train"/"Entities" block (rows 24-65) with the restructured block (rows 24-68) ---
# Delete old rows entirely (removes stale formatting/heights) then insert fresh blank rows
$ws.Rows("24:65").Delete()
$ws.Rows("24:68").Insert()

$ws.Cells.Item(24, 1).Value = 'Train'
$ws.Cells.Item(24, 2).Value = 'Edit Turn'
$ws.Cells.Item(24, 3).Value = 'Train Dialog'
$ws.Cells.Item(24, 4).Value = 'Selecting a turn causes special buttons to show up'
$ws.Cells.Item(24, 5).Value = 'VerifyEditTrainingControlsAndLabels'

$ws.Cells.Item(25, 1).Value = 'Train'
$ws.Cells.Item(25, 2).Value = 'Edit Turn'
$ws.Cells.Item(25, 3).Value = 'Train Dialog'
$ws.Cells.Item(25, 4).Value = 'Delete button shows for all EXCEPT the first turn'
$ws.Cells.Item(25, 5).Value = 'VerifyEditTrainingControlsAndLabels'

$ws.Cells.Item(26, 1).Value = 'Train'
$ws.Cells.Item(26, 2).Value = 'Edit Turn'
$ws.Cells.Item(26, 3).Value = 'Train Dialog'
$ws.Cells.Item(26, 4).Value = 'Add score and add input buttons show for all turns'
$ws.Cells.Item(26, 5).Value = 'VerifyEditTrainingControlsAndLabels'

$ws.Cells.Item(27, 1).Value = 'Train'
$ws.Cells.Item(27, 2).Value = 'Edit Turn'
$ws.Cells.Item(27, 3).Value = 'Train Dialog'
$ws.Cells.Item(27, 4).Value = 'Branching button shows up for all user turns but NOT for Bot turns'
$ws.Cells.Item(27, 5).Value = 'VerifyEditTrainingControlsAndLabels'

$ws.Cells.Item(28, 1).Value = 'Train'
$ws.Cells.Item(28, 2).Value = 'Edit Turn'
$ws.Cells.Item(28, 3).Value = 'Train Dialog'
$ws.Cells.Item(28, 4).Value = 'Not selecting a turn and no special buttons show up'
$ws.Cells.Item(28, 5).Value = 'VerifyEditTrainingControlsAndLabels'

$ws.Cells.Item(29, 1).Value = 'Train'
$ws.Cells.Item(29, 2).Value = 'Edit Turn'
$ws.Cells.Item(29, 3).Value = 'Train Dialog'
$ws.Cells.Item(29, 4).Value = 'Selecting a user turn causes "Entity Detection" UI Elements to show up'

$ws.Cells.Item(30, 1).Value = 'Train'
$ws.Cells.Item(30, 2).Value = 'Edit Turn'
$ws.Cells.Item(30, 3).Value = 'Train Dialog'
$ws.Cells.Item(30, 4).Value = 'Selecting a Bot turn causes Action Selection UI Elements to show up'

$ws.Cells.Item(33, 1).Value = 'Train'
$ws.Cells.Item(33, 2).Value = 'Branching'
$ws.Cells.Item(33, 3).Value = 'Train Dialog'
$ws.Cells.Item(33, 4).Value = 'Buttons should remain "Close" and "Delete" as long as nothing changes'
$ws.Cells.Item(33, 5).Value = 'VerifyEditTrainingControlsAndLabels'

$ws.Cells.Item(34, 1).Value = 'Train'
$ws.Cells.Item(34, 2).Value = 'Branching'
$ws.Cells.Item(34, 3).Value = 'Train Dialog'
$ws.Cells.Item(34, 4).Value = '"Close" and "Delete" buttons should change to "Save Branch" and "Abandon Branch" after branching'
$ws.Cells.Item(34, 5).Value = 'VerifyEditTrainingControlsAndLabels'
$ws.Rows.Item(34).RowHeight = 28.8

$ws.Cells.Item(35, 1).Value = 'Train'
$ws.Cells.Item(35, 2).Value = 'Branching'
$ws.Cells.Item(35, 3).Value = 'Train Dialog'
$ws.Cells.Item(35, 4).Value = 'Create new branch should create a new branch training that changes the selected turn to "Different User Input" that was entered'
$ws.Cells.Item(35, 5).Value = 'Branching'
$ws.Rows.Item(35).RowHeight = 28.8

$ws.Cells.Item(36, 1).Value = 'Train'
$ws.Cells.Item(36, 2).Value = 'Branching'
$ws.Cells.Item(36, 3).Value = 'Train Dialog'
$ws.Cells.Item(36, 4).Value = 'Create new branch should create a new branch training that cuts off everything below the selected turn'
$ws.Cells.Item(36, 5).Value = 'Branching'
$ws.Rows.Item(36).RowHeight = 28.8

$ws.Cells.Item(37, 1).Value = 'Train'
$ws.Cells.Item(37, 2).Value = 'Branching'
$ws.Cells.Item(37, 3).Value = 'Train Dialog'
$ws.Cells.Item(37, 4).Value = 'After new branch has been created, all edit controls in the chat pane should disappear'
$ws.Cells.Item(37, 5).Value = 'Branching'

$ws.Cells.Item(38, 1).Value = 'Train'
$ws.Cells.Item(38, 2).Value = 'Branching'
$ws.Cells.Item(38, 3).Value = 'Train Dialog'
$ws.Cells.Item(38, 4).Value = 'Continued training of a new branch should work and adds new dialog to the training that is persisted'
$ws.Rows.Item(38).RowHeight = 28.8

$ws.Cells.Item(39, 1).Value = 'Train'
$ws.Cells.Item(39, 2).Value = 'Branching'
$ws.Cells.Item(39, 3).Value = 'Train Dialog'
$ws.Cells.Item(39, 4).Value = '"Save Branch" button should save the branch and leave the original branch in the grid'
$ws.Cells.Item(39, 5).Value = 'Branching'

$ws.Cells.Item(40, 1).Value = 'Train'
$ws.Cells.Item(40, 2).Value = 'Branching'
$ws.Cells.Item(40, 3).Value = 'Train Dialog'
$ws.Cells.Item(40, 4).Value = 'After branch is saved the original training should remain unchanged'
$ws.Cells.Item(40, 5).Value = 'Branching'

$ws.Cells.Item(41, 1).Value = 'Train'
$ws.Cells.Item(41, 2).Value = 'Branching'
$ws.Cells.Item(41, 3).Value = 'Train Dialog'
$ws.Cells.Item(41, 4).Value = 'After branch is abandonded the original training should remain unchanged'
$ws.Cells.Item(41, 5).Value = 'VerifyEditTrainingControlsAndLabels'

$ws.Cells.Item(42, 1).Value = 'Train'
$ws.Cells.Item(42, 2).Value = 'Branching'
$ws.Cells.Item(42, 3).Value = 'Train Dialog'
$ws.Cells.Item(42, 4).Value = '"Abandon Branch" button should leave the original branch in the grid'
$ws.Cells.Item(42, 5).Value = 'VerifyEditTrainingControlsAndLabels'

$ws.Cells.Item(43, 1).Value = 'Train'
$ws.Cells.Item(43, 2).Value = 'Branching'
$ws.Cells.Item(43, 3).Value = 'Train Dialog'
$ws.Cells.Item(43, 4).Value = '"Abandon Branch" button should not persist the new Train Dialog, it should NOT show up in the grid'
$ws.Rows.Item(43).RowHeight = 28.8

$ws.Cells.Item(44, 1).Value = 'Train'
$ws.Cells.Item(44, 2).Value = 'Branching'
$ws.Cells.Item(44, 3).Value = 'Train Dialog'
$ws.Cells.Item(44, 4).Value = 'Attempts to branch above training errors should succeed'
$ws.Cells.Item(44, 6).Value = 'Training errors are due to Deleted Entities or Actions that a training depended on'

$ws.Cells.Item(45, 1).Value = 'Train'
$ws.Cells.Item(45, 2).Value = 'Branching'
$ws.Cells.Item(45, 3).Value = 'Train Dialog'
$ws.Cells.Item(45, 4).Value = 'Attempts to branch at or below training errors should fail'

$ws.Cells.Item(46, 1).Value = 'Train'
$ws.Cells.Item(46, 2).Value = 'Branching'
$ws.Cells.Item(46, 3).Value = 'Train Dialog'
$ws.Cells.Item(46, 4).Value = 'TODO: Branching + Edit how do they mix? Edit first does not allow branching.'
$ws.Cells.Item(46, 4).Font.Color = 255
$ws.Cells.Item(46, 4).WrapText = $true

$ws.Cells.Item(48, 1).Value = 'Train'
$ws.Cells.Item(48, 2).Value = 'Label Entities'
$ws.Cells.Item(48, 3).Value = 'Train Dialog'

$ws.Cells.Item(50, 1).Value = 'Entities'
$ws.Cells.Item(50, 2).Value = 'Multi-Value'
$ws.Cells.Item(50, 3).Value = 'Entities Dialog'
$ws.Cells.Item(50, 4).Value = 'Setting is persisted, shows checked after save then edit'

$ws.Cells.Item(51, 1).Value = 'Entities'
$ws.Cells.Item(51, 2).Value = 'Multi-Value'
$ws.Cells.Item(51, 3).Value = 'Entities Grid'
$ws.Cells.Item(51, 4).Value = 'Checkbox is set in grid'

$ws.Cells.Item(52, 1).Value = 'Entities'
$ws.Cells.Item(52, 2).Value = 'Multi-Value'
$ws.Cells.Item(52, 3).Value = 'Train Dialog - Memory Panel'
$ws.Cells.Item(52, 4).Value = 'Entities accumulate values in training memory'

$ws.Cells.Item(53, 1).Value = 'Entities'
$ws.Cells.Item(53, 2).Value = 'Multi-Value'
$ws.Cells.Item(53, 3).Value = 'Train Dialog - Memory Panel'
$ws.Cells.Item(53, 4).Value = 'Checkbox is set in grid'

$ws.Cells.Item(54, 1).Value = 'Entities'
$ws.Cells.Item(54, 2).Value = 'Multi-Value'
$ws.Cells.Item(54, 3).Value = 'Train Dialog - Chat Panel'
$ws.Cells.Item(54, 4).Value = 'All values show up as a list when displayed in training webchat'

$ws.Cells.Item(55, 1).Value = 'Entities'
$ws.Cells.Item(55, 2).Value = 'Multi-Value'
$ws.Cells.Item(55, 3).Value = 'Log Dialog - Chat Panel'
$ws.Cells.Item(55, 4).Value = 'All values show up as a list when displayed in training webchat'

$ws.Cells.Item(57, 1).Value = 'Entities'
$ws.Cells.Item(57, 2).Value = 'Negatable'
$ws.Cells.Item(57, 3).Value = 'Entities Dialog'
$ws.Cells.Item(57, 4).Value = 'Setting is persisted, shows checked after save then edit'

$ws.Cells.Item(58, 1).Value = 'Entities'
$ws.Cells.Item(58, 2).Value = 'Negatable'
$ws.Cells.Item(58, 3).Value = 'Entities Grid'
$ws.Cells.Item(58, 4).Value = 'Checkbox is set in grid'

$ws.Cells.Item(59, 1).Value = 'Entities'
$ws.Cells.Item(59, 2).Value = 'Negatable'
$ws.Cells.Item(59, 3).Value = 'Train Dialog - Memory Panel'
$ws.Cells.Item(59, 4).Value = 'Entities remove values in training memory'

$ws.Cells.Item(60, 1).Value = 'Entities'
$ws.Cells.Item(60, 2).Value = 'Negatable'
$ws.Cells.Item(60, 3).Value = 'Train Dialog - Memory Panel'
$ws.Cells.Item(60, 4).Value = 'Checkbox is set in grid'

$ws.Cells.Item(61, 1).Value = 'Entities'
$ws.Cells.Item(61, 2).Value = 'Negatable'
$ws.Cells.Item(61, 3).Value = 'Train Dialog - Chat Panel'
$ws.Cells.Item(61, 4).Value = 'All values except for negated one shows up as a list when displayed in training webchat'

$ws.Cells.Item(62, 1).Value = 'Entities'
$ws.Cells.Item(62, 2).Value = 'Negatable'
$ws.Cells.Item(62, 3).Value = 'Log Dialog - Chat Panel'
$ws.Cells.Item(62, 4).Value = 'All values except for negated one shows up as a list when displayed in training webchat'

$ws.Cells.Item(63, 1).Value = 'Entities'

$ws.Cells.Item(64, 1).Value = 'Entities'

$ws.Cells.Item(65, 1).Value = 'Entities'

$ws.Cells.Item(66, 1).Value = 'Entities'

$ws.Cells.Item(67, 1).Value = 'Entities'

$ws.Cells.Item(68, 1).Value = 'Entities'

# --- Resize Table1 to the new range ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F68"))

# --- Update selection to match target view ---
$ws.Range("D30").Select() | Out-Null